$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to write a value while preserving it as text (matches the
# original inline-string cell type used throughout column D, etc.)
function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# --- Simple price (column D) updates -------------------------------------
$priceUpdates = @{
    2  = "247.82"
    3  = "22.78"
    4  = "5.296"
    5  = "0.05723"
    6  = "3.426"
    7  = "0.8097"
    8  = "0.8734"
    9  = "0.1432"
    10 = "0.07411"
    11 = "0.03057"
    12 = "0.03114"
    13 = "0.09387"
    14 = "3.880"
    15 = "0.001589"
    16 = "0.04800"
    25 = "0.3278"
    26 = "0.1330"
    40 = "0.03938"
    41 = "0.006742"
    43 = "0.002220"
    44 = "0.007983"
    45 = "0.00005607"
    47 = "0.6001"
    48 = "0.1801"
    49 = "0.00002100"
}

foreach ($row in $priceUpdates.Keys) {
    Set-TextValue $ws.Cells.Item($row, 4) $priceUpdates[$row]
}

# --- Rows 17-24: "One" moved up to rank 17 (ahead of TigerCash), pushing
#     the previously ranked coins (TigerCash..BTSEToken) down by one row,
#     each refreshed with a new price and rank label in column E. --------
$rowData = @{
    17 = @{ B = "One";         C = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one";          D = "0.0005851"; E = "16OneONE" }
    18 = @{ B = "TigerCash";   C = "https://coinranking.com/coin/6hIn06L2+tigercash-tch";          D = "0.006147";  E = "17TigerCashTCH" }
    19 = @{ B = "HotbitToken"; C = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb";    D = "0.005168";  E = "18HotbitTokenHTB" }
    20 = @{ B = "BitKan";      C = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan";        D = "0.0009983"; E = "19BitKanKAN" }
    21 = @{ B = "NitroEx";     C = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx";         D = "0.0001500"; E = "20NitroExNTX" }
    22 = @{ B = "LEO";         C = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo";            D = "3.738";     E = "21LEOLEO" }
    23 = @{ B = "KuCoinToken"; C = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs";   D = "6.334";     E = "22KuCoinTokenKCS" }
    24 = @{ B = "BTSEToken";   C = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse";     D = "2.195";     E = "23BTSETokenBTSE" }
}

foreach ($row in $rowData.Keys) {
    $data = $rowData[$row]
    Set-TextValue $ws.Cells.Item($row, 2) $data.B
    Set-TextValue $ws.Cells.Item($row, 3) $data.C
    Set-TextValue $ws.Cells.Item($row, 4) $data.D
    Set-TextValue $ws.Cells.Item($row, 5) $data.E
}

# --- Row 47 column E gains a "Bestin24h" suffix ---------------------------
Set-TextValue $ws.Cells.Item(47, 5) "46CoinbaseStockTokenCOINBestin24h"
